# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2026-01-26 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-01-27 Tuesday", 2)

# Update the multiplication problems in the table. Target cells explicitly
# by (row, column) so the two duplicate "52×22=" cells each get the correct,
# independent replacement value.
$tbl = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="56×56="},
    @{Row=1;  Col=2; New="31×82="},
    @{Row=1;  Col=3; New="33×43="},
    @{Row=1;  Col=4; New="17×32="},
    @{Row=1;  Col=5; New="33×97="},

    @{Row=5;  Col=1; New="27×49="},
    @{Row=5;  Col=2; New="41×85="},
    @{Row=5;  Col=3; New="99×85="},
    @{Row=5;  Col=4; New="98×31="},
    @{Row=5;  Col=5; New="18×84="},

    @{Row=10; Col=1; New="58×77="},
    @{Row=10; Col=2; New="71×69="},
    @{Row=10; Col=3; New="99×53="},
    @{Row=10; Col=4; New="19×41="},
    @{Row=10; Col=5; New="63×77="},

    @{Row=15; Col=1; New="43×50="},
    @{Row=15; Col=2; New="31×41="},
    @{Row=15; Col=3; New="56×35="},
    @{Row=15; Col=4; New="89×78="},
    @{Row=15; Col=5; New="92×92="},

    @{Row=20; Col=1; New="75×60="},
    @{Row=20; Col=2; New="71×19="},
    @{Row=20; Col=3; New="79×23="},
    @{Row=20; Col=4; New="99×88="},
    @{Row=20; Col=5; New="26×98="}
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $cellRange = $cell.Range
    $cellRange.End = $cellRange.End - 1
    $cellRange.Text = $u.New
}
